# Add data for 2022-03-14: update the "as of" date in the sheet name and
# header label, and bump the March (row 4) and Total (row 14) counts in the
# "Total" / through-date column (I).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new as-of date.
$ws.Name = "Through 2022-03-06"

# Update the header label in I1 (shared string "2022 (through 03-05)" -> "... 03-06").
$ws.Range("I1").Value = "2022 (through 03-06)"

# Update the March row's running total and the grand Total row.
$ws.Range("I4").Value = 33
$ws.Range("I14").Value = 334
